# Edit script for Cn_fullerenes.xlsx
# Fills in Point Group (column B) and Energy (column C) data for the
# C46 and C48 sheets, which previously only had the index column (A)
# populated. Also updates sheet/view selection state to match the
# author's last-saved position, and propagates the new "C2h" point
# group value (used by one of the new C48 rows) into the shared
# string table.

$wb = $excel.ActiveWorkbook

# ---- C46 sheet: populate columns B (Point Group) and C (Energy) ----
$c46data = New-Object 'object[,]' 116,2
$c46data[0,0] = "C2"
$c46data[0,1] = 30.405
$c46data[1,0] = "Cs"
$c46data[1,1] = 27.337
$c46data[2,0] = "C1"
$c46data[2,1] = 28.111
$c46data[3,0] = "C1"
$c46data[3,1] = 27.805
$c46data[4,0] = "C1"
$c46data[4,1] = 27.129
$c46data[5,0] = "C1"
$c46data[5,1] = 28.221
$c46data[6,0] = "Cs"
$c46data[6,1] = 27.585
$c46data[7,0] = "Cs"
$c46data[7,1] = 26.77
$c46data[8,0] = "C2"
$c46data[8,1] = 28.473
$c46data[9,0] = "Cs"
$c46data[9,1] = 27.819
$c46data[10,0] = "Cs"
$c46data[10,1] = 27.401
$c46data[11,0] = "C2"
$c46data[11,1] = 26.825
$c46data[12,0] = "Cs"
$c46data[12,1] = 25.99
$c46data[13,0] = "C1"
$c46data[13,1] = 26.263
$c46data[14,0] = "C1"
$c46data[14,1] = 28.076
$c46data[15,0] = "C1"
$c46data[15,1] = 26.843
$c46data[16,0] = "C1"
$c46data[16,1] = 26.501
$c46data[17,0] = "C1"
$c46data[17,1] = 26.377
$c46data[18,0] = "C1"
$c46data[18,1] = 26.621
$c46data[19,0] = "C2"
$c46data[19,1] = 28.709
$c46data[20,0] = "C1"
$c46data[20,1] = 27.152
$c46data[21,0] = "C2"
$c46data[21,1] = 26.061
$c46data[22,0] = "C1"
$c46data[22,1] = 25.562
$c46data[23,0] = "C1"
$c46data[23,1] = 25.631
$c46data[24,0] = "C1"
$c46data[24,1] = 25.64
$c46data[25,0] = "C1"
$c46data[25,1] = 26.607
$c46data[26,0] = "C1"
$c46data[26,1] = 26.224
$c46data[27,0] = "Cs"
$c46data[27,1] = 25.888
$c46data[28,0] = "C1"
$c46data[28,1] = 25.451
$c46data[29,0] = "C1"
$c46data[29,1] = 25.247
$c46data[30,0] = "C1"
$c46data[30,1] = 27.241
$c46data[31,0] = "C2"
$c46data[31,1] = 25.439
$c46data[32,0] = "Cs"
$c46data[32,1] = 27.644
$c46data[33,0] = "C1"
$c46data[33,1] = 25.632
$c46data[34,0] = "C1"
$c46data[34,1] = 24.951
$c46data[35,0] = "C1"
$c46data[35,1] = 24.757
$c46data[36,0] = "C1"
$c46data[36,1] = 26.261
$c46data[37,0] = "Cs"
$c46data[37,1] = 25.907
$c46data[38,0] = "C2v"
$c46data[38,1] = 25.323
$c46data[39,0] = "Cs"
$c46data[39,1] = 24.492
$c46data[40,0] = "Cs"
$c46data[40,1] = 25.471
$c46data[41,0] = "C2v"
$c46data[41,1] = 27.582
$c46data[42,0] = "C2"
$c46data[42,1] = 25.589
$c46data[43,0] = "C1"
$c46data[43,1] = 25.351
$c46data[44,0] = "C1"
$c46data[44,1] = 26.618
$c46data[45,0] = "C1"
$c46data[45,1] = 25.132
$c46data[46,0] = "C2"
$c46data[46,1] = 26.054
$c46data[47,0] = "C1"
$c46data[47,1] = 26.898
$c46data[48,0] = "C2"
$c46data[48,1] = 26.299
$c46data[49,0] = "C1"
$c46data[49,1] = 25.054
$c46data[50,0] = "C1"
$c46data[50,1] = 25.219
$c46data[51,0] = "C1"
$c46data[51,1] = 25.719
$c46data[52,0] = "C2"
$c46data[52,1] = 26.177
$c46data[53,0] = "C2"
$c46data[53,1] = 24.938
$c46data[54,0] = "C1"
$c46data[54,1] = 24.869
$c46data[55,0] = "C1"
$c46data[55,1] = 25.737
$c46data[56,0] = "Cs"
$c46data[56,1] = 26.421
$c46data[57,0] = "C1"
$c46data[57,1] = 25.544
$c46data[58,0] = "C1"
$c46data[58,1] = 25.024
$c46data[59,0] = "C1"
$c46data[59,1] = 24.257
$c46data[60,0] = "C1"
$c46data[60,1] = 24.272
$c46data[61,0] = "C1"
$c46data[61,1] = 24.946
$c46data[62,0] = "C1"
$c46data[62,1] = 24.378
$c46data[63,0] = "C1"
$c46data[63,1] = 25.331
$c46data[64,0] = "Cs"
$c46data[64,1] = 24.613
$c46data[65,0] = "C2"
$c46data[65,1] = 24.204
$c46data[66,0] = "C1"
$c46data[66,1] = 23.851
$c46data[67,0] = "C1"
$c46data[67,1] = 24.521
$c46data[68,0] = "C1"
$c46data[68,1] = 24.186
$c46data[69,0] = "C1"
$c46data[69,1] = 24.039
$c46data[70,0] = "C1"
$c46data[70,1] = 24.983
$c46data[71,0] = "C1"
$c46data[71,1] = 24.167
$c46data[72,0] = "C1"
$c46data[72,1] = 24.637
$c46data[73,0] = "C1"
$c46data[73,1] = 24.188
$c46data[74,0] = "C1"
$c46data[74,1] = 24.612
$c46data[75,0] = "C1"
$c46data[75,1] = 24.884
$c46data[76,0] = "C2"
$c46data[76,1] = 25.388
$c46data[77,0] = "C1"
$c46data[77,1] = 24.596
$c46data[78,0] = "C1"
$c46data[78,1] = 24.275
$c46data[79,0] = "C1"
$c46data[79,1] = 24.963
$c46data[80,0] = "C1"
$c46data[80,1] = 24.523
$c46data[81,0] = "C1"
$c46data[81,1] = 24.558
$c46data[82,0] = "Cs"
$c46data[82,1] = 25.265
$c46data[83,0] = "C2"
$c46data[83,1] = 24.33
$c46data[84,0] = "C1"
$c46data[84,1] = 23.608
$c46data[85,0] = "C1"
$c46data[85,1] = 23.325
$c46data[86,0] = "C1"
$c46data[86,1] = 23.562
$c46data[87,0] = "C1"
$c46data[87,1] = 23.282
$c46data[88,0] = "Cs"
$c46data[88,1] = 24.467
$c46data[89,0] = "C1"
$c46data[89,1] = 23.253
$c46data[90,0] = "C2v"
$c46data[90,1] = 24.685
$c46data[91,0] = "C2v"
$c46data[91,1] = 23.934
$c46data[92,0] = "C1"
$c46data[92,1] = 23.589
$c46data[93,0] = "C3"
$c46data[93,1] = 23.011
$c46data[94,0] = "C2"
$c46data[94,1] = 23.55
$c46data[95,0] = "C2"
$c46data[95,1] = 24.017
$c46data[96,0] = "C2"
$c46data[96,1] = 24.433
$c46data[97,0] = "C1"
$c46data[97,1] = 23.813
$c46data[98,0] = "Cs"
$c46data[98,1] = 23.465
$c46data[99,0] = "C1"
$c46data[99,1] = 23.779
$c46data[100,0] = "C1"
$c46data[100,1] = 23.419
$c46data[101,0] = "C1"
$c46data[101,1] = 23.58
$c46data[102,0] = "C1"
$c46data[102,1] = 23.069
$c46data[103,0] = "C2"
$c46data[103,1] = 23.245
$c46data[104,0] = "C1"
$c46data[104,1] = 23.855
$c46data[105,0] = "Cs"
$c46data[105,1] = 23.274
$c46data[106,0] = "Cs"
$c46data[106,1] = 23.273
$c46data[107,0] = "Cs"
$c46data[107,1] = 22.958
$c46data[108,0] = "C2"
$c46data[108,1] = 22.769
$c46data[109,0] = "C1"
$c46data[109,1] = 23.101
$c46data[110,0] = "C1"
$c46data[110,1] = 23.674
$c46data[111,0] = "C2"
$c46data[111,1] = 23.576
$c46data[112,0] = "C2"
$c46data[112,1] = 25.297
$c46data[113,0] = "C1"
$c46data[113,1] = 22.671
$c46data[114,0] = "C3"
$c46data[114,1] = 23.093
$c46data[115,0] = "C2"
$c46data[115,1] = 22.705
$wsC46 = $wb.Worksheets.Item("C46")
$wsC46.Range("B2:C117").Value = $c46data

# ---- C48 sheet: populate columns B (Point Group) and C (Energy) ----
$c48data = New-Object 'object[,]' 120,2
$c48data[0,0] = "C2"
$c48data[0,1] = 30.885
$c48data[1,0] = "D2"
$c48data[1,1] = 36.058
$c48data[2,0] = "C1"
$c48data[2,1] = 28.913
$c48data[3,0] = "Cs"
$c48data[3,1] = 28.725
$c48data[4,0] = "C2"
$c48data[4,1] = 31.082
$c48data[5,0] = "C1"
$c48data[5,1] = 29.523
$c48data[6,0] = "C1"
$c48data[6,1] = 28.404
$c48data[7,0] = "C1"
$c48data[7,1] = 27.765
$c48data[8,0] = "C1"
$c48data[8,1] = 29.493
$c48data[9,0] = "C1"
$c48data[9,1] = 27.401
$c48data[10,0] = "C1"
$c48data[10,1] = 27.834
$c48data[11,0] = "C1"
$c48data[11,1] = 29.653
$c48data[12,0] = "C1"
$c48data[12,1] = 28.115
$c48data[13,0] = "C2"
$c48data[13,1] = 29.372
$c48data[14,0] = "D2h"
$c48data[14,1] = 29.264
$c48data[15,0] = "D2"
$c48data[15,1] = 28.861
$c48data[16,0] = "C2v"
$c48data[16,1] = 27.597
$c48data[17,0] = "C1"
$c48data[17,1] = 27.262
$c48data[18,0] = "C1"
$c48data[18,1] = 26.626
$c48data[19,0] = "C1"
$c48data[19,1] = 27.66
$c48data[20,0] = "C1"
$c48data[20,1] = 27.155
$c48data[21,0] = "C1"
$c48data[21,1] = 26.229
$c48data[22,0] = "C1"
$c48data[22,1] = 27.337
$c48data[23,0] = "C2"
$c48data[23,1] = 27.637
$c48data[24,0] = "C1"
$c48data[24,1] = 27.275
$c48data[25,0] = "C1"
$c48data[25,1] = 28.008
$c48data[26,0] = "C2"
$c48data[26,1] = 28.479
$c48data[27,0] = "C1"
$c48data[27,1] = 26.279
$c48data[28,0] = "C1"
$c48data[28,1] = 27.085
$c48data[29,0] = "C1"
$c48data[29,1] = 28.325
$c48data[30,0] = "Cs"
$c48data[30,1] = 28.005
$c48data[31,0] = "C2"
$c48data[31,1] = 26.419
$c48data[32,0] = "C1"
$c48data[32,1] = 25.713
$c48data[33,0] = "C1"
$c48data[33,1] = 25.984
$c48data[34,0] = "C1"
$c48data[34,1] = 27.251
$c48data[35,0] = "C1"
$c48data[35,1] = 26.464
$c48data[36,0] = "C2"
$c48data[36,1] = 26.189
$c48data[37,0] = "C1"
$c48data[37,1] = 25.618
$c48data[38,0] = "Cs"
$c48data[38,1] = 25.491
$c48data[39,0] = "C2"
$c48data[39,1] = 26.118
$c48data[40,0] = "D2h"
$c48data[40,1] = 26.017
$c48data[41,0] = "C1"
$c48data[41,1] = 25.443
$c48data[42,0] = "C2"
$c48data[42,1] = 26.179
$c48data[43,0] = "C1"
$c48data[43,1] = 25.464
$c48data[44,0] = "C2"
$c48data[44,1] = 25.72
$c48data[45,0] = "C2"
$c48data[45,1] = 28.924
$c48data[46,0] = "C1"
$c48data[46,1] = 25.515
$c48data[47,0] = "C1"
$c48data[47,1] = 26.744
$c48data[48,0] = "C1"
$c48data[48,1] = 29.797
$c48data[49,0] = "C1"
$c48data[49,1] = 25.872
$c48data[50,0] = "C1"
$c48data[50,1] = 25.167
$c48data[51,0] = "C1"
$c48data[51,1] = 26.139
$c48data[52,0] = "C1"
$c48data[52,1] = 25.501
$c48data[53,0] = "C1"
$c48data[53,1] = 26.065
$c48data[54,0] = "C1"
$c48data[54,1] = 25.991
$c48data[55,0] = "C2v"
$c48data[55,1] = 26.45
$c48data[56,0] = "C1"
$c48data[56,1] = 27.946
$c48data[57,0] = "C2"
$c48data[57,1] = 25.791
$c48data[58,0] = "C2"
$c48data[58,1] = 28.118
$c48data[59,0] = "C1"
$c48data[59,1] = 26.821
$c48data[60,0] = "C2"
$c48data[60,1] = 25.325
$c48data[61,0] = "Cs"
$c48data[61,1] = 25.195
$c48data[62,0] = "C2"
$c48data[62,1] = 28.206
$c48data[63,0] = "C2"
$c48data[63,1] = 26.554
$c48data[64,0] = "C1"
$c48data[64,1] = 25.931
$c48data[65,0] = "C1"
$c48data[65,1] = 26.278
$c48data[66,0] = "C1"
$c48data[66,1] = 25.344
$c48data[67,0] = "C2"
$c48data[67,1] = 26.913
$c48data[68,0] = "C1"
$c48data[68,1] = 27.326
$c48data[69,0] = "C2"
$c48data[69,1] = 27.083
$c48data[70,0] = "C1"
$c48data[70,1] = 24.95
$c48data[71,0] = "C1"
$c48data[71,1] = 25.032
$c48data[72,0] = "C1"
$c48data[72,1] = 24.87
$c48data[73,0] = "Cs"
$c48data[73,1] = 25.326
$c48data[74,0] = "Cs"
$c48data[74,1] = 25.979
$c48data[75,0] = "C2"
$c48data[75,1] = 26.552
$c48data[76,0] = "C1"
$c48data[76,1] = 25.225
$c48data[77,0] = "C2"
$c48data[77,1] = 25.439
$c48data[78,0] = "C1"
$c48data[78,1] = 25.984
$c48data[79,0] = "C2h"
$c48data[79,1] = 25.272
$c48data[80,0] = "C2"
$c48data[80,1] = 25.123
$c48data[81,0] = "C2"
$c48data[81,1] = 24.933
$c48data[82,0] = "C2"
$c48data[82,1] = 24.763
$c48data[83,0] = "C2"
$c48data[83,1] = 24.755
$c48data[84,0] = "C1"
$c48data[84,1] = 24.174
$c48data[85,0] = "C1"
$c48data[85,1] = 26.222
$c48data[86,0] = "C1"
$c48data[86,1] = 24.266
$c48data[87,0] = "C1"
$c48data[87,1] = 24.738
$c48data[88,0] = "Cs"
$c48data[88,1] = 26.535
$c48data[89,0] = "C1"
$c48data[89,1] = 24.378
$c48data[90,0] = "C1"
$c48data[90,1] = 23.956
$c48data[91,0] = "C1"
$c48data[91,1] = 24.514
$c48data[92,0] = "C1"
$c48data[92,1] = 24.209
$c48data[93,0] = "C1"
$c48data[93,1] = 24.683
$c48data[94,0] = "C2"
$c48data[94,1] = 24.86
$c48data[95,0] = "Cs"
$c48data[95,1] = 26.325
$c48data[96,0] = "C2"
$c48data[96,1] = 25.655
$c48data[97,0] = "C1"
$c48data[97,1] = 25.613
$c48data[98,0] = "C1"
$c48data[98,1] = 25.682
$c48data[99,0] = "C1"
$c48data[99,1] = 25.171
$c48data[100,0] = "C1"
$c48data[100,1] = 26.326
$c48data[101,0] = "C1"
$c48data[101,1] = 24.733
$c48data[102,0] = "C1"
$c48data[102,1] = 24.322
$c48data[103,0] = "C1"
$c48data[103,1] = 24.324
$c48data[104,0] = "C1"
$c48data[104,1] = 26.312
$c48data[105,0] = "C1"
$c48data[105,1] = 25.133
$c48data[106,0] = "C2"
$c48data[106,1] = 25.254
$c48data[107,0] = "C1"
$c48data[107,1] = 24.744
$c48data[108,0] = "C1"
$c48data[108,1] = 25.47
$c48data[109,0] = "C1"
$c48data[109,1] = 24.589
$c48data[110,0] = "C1"
$c48data[110,1] = 25.282
$c48data[111,0] = "C1"
$c48data[111,1] = 25.366
$c48data[112,0] = "C1"
$c48data[112,1] = 24.524
$c48data[113,0] = "C1"
$c48data[113,1] = 25.837
$c48data[114,0] = "C2"
$c48data[114,1] = 26.575
$c48data[115,0] = "C1"
$c48data[115,1] = 24.911
$c48data[116,0] = "C1"
$c48data[116,1] = 25.315
$c48data[117,0] = "C1"
$c48data[117,1] = 24.916
$c48data[118,0] = "C1"
$c48data[118,1] = 26.437
$c48data[119,0] = "C1"
$c48data[119,1] = 26.51
$wsC48 = $wb.Worksheets.Item("C48")
$wsC48.Range("B2:C121").Value = $c48data

# ---- View/selection state ----
# Final workbook has C48 as the active (selected) tab, with the C46
# sheet scrolled to show its last row and C48 scrolled roughly to its
# last populated row.
[void]$wsC46.Activate()
[void]$wsC46.Range("C118").Select()

[void]$wsC48.Activate()
[void]$wsC48.Range("B122").Select()

$wsC96 = $wb.Worksheets.Item("C96")
[void]$wsC96.Activate()
[void]$wsC96.Range("I23").Select()

# Re-activate C48 last so it ends up the active sheet/tab.
[void]$wsC48.Activate()
